$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 5

$ws.Range("D6").Select()
